$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly by writing data rows 2-13 (replacing old rows 2-9,
# the workbook originally had 8 data rows; now it has 12).

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Has2"
$ws.Range("C2").Value = "Cd44"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 23.544642
$ws.Range("H2").Value = 70.633926
$ws.Range("I2").Value = 0.8736569986425735
$ws.Range("J2").Value = 0.8736569986425735
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 239.0839323333333
$ws.Range("N2").Value = 717.251797
$ws.Range("O2").Value = 0.4086975387666237
$ws.Range("P2").Value = 0.4086975387666237
$ws.Range("Q2").Value = 5629.145594740558
$ws.Range("R2").Value = 50662.31035266502
$ws.Range("S2").Value = 0.3570614650714553
$ws.Range("T2").Value = 0.3570614650714553

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Has2"
$ws.Range("C3").Value = "Cd44"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 23.544642
$ws.Range("H3").Value = 70.633926
$ws.Range("I3").Value = 0.8736569986425735
$ws.Range("J3").Value = 0.8736569986425735
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 117.0512696666667
$ws.Range("N3").Value = 351.153809
$ws.Range("O3").Value = 0.2000910950200451
$ws.Range("P3").Value = 0.2000910950200451
$ws.Range("Q3").Value = 2755.930239947126
$ws.Range("R3").Value = 24803.37215952414
$ws.Range("S3").Value = 0.1748109855303186
$ws.Range("T3").Value = 0.1748109855303186

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Has2"
$ws.Range("C4").Value = "Cd44"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 23.544642
$ws.Range("H4").Value = 70.633926
$ws.Range("I4").Value = 0.8736569986425735
$ws.Range("J4").Value = 0.8736569986425735
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 171.15883
$ws.Range("N4").Value = 513.47649
$ws.Range("O4").Value = 0.2925842480357353
$ws.Range("P4").Value = 0.2925842480357353
$ws.Range("Q4").Value = 4029.87337748886
$ws.Range("R4").Value = 36268.86039739975
$ws.Range("S4").Value = 0.2556182759889948
$ws.Range("T4").Value = 0.2556182759889948

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Has2"
$ws.Range("C5").Value = "Cd44"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.544642
$ws.Range("H5").Value = 70.633926
$ws.Range("I5").Value = 0.8736569986425735
$ws.Range("J5").Value = 0.8736569986425735
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 57.695868
$ws.Range("N5").Value = 173.087604
$ws.Range("O5").Value = 0.09862711817759588
$ws.Range("P5").Value = 0.09862711817759588
$ws.Range("Q5").Value = 1358.428556939256
$ws.Range("R5").Value = 12225.8570124533
$ws.Range("S5").Value = 0.08616627205180483
$ws.Range("T5").Value = 0.08616627205180483

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Has2"
$ws.Range("C6").Value = "Cd44"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.041467
$ws.Range("H6").Value = 0.124401
$ws.Range("I6").Value = 0.001538691255645832
$ws.Range("J6").Value = 0.001538691255645832
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 239.0839323333333
$ws.Range("N6").Value = 717.251797
$ws.Range("O6").Value = 0.4086975387666237
$ws.Range("P6").Value = 0.4086975387666237
$ws.Range("Q6").Value = 9.914093422066333
$ws.Range("R6").Value = 89.226840798597
$ws.Range("S6").Value = 0.0006288593291041772
$ws.Range("T6").Value = 0.0006288593291041772

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Has2"
$ws.Range("C7").Value = "Cd44"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.041467
$ws.Range("H7").Value = 0.124401
$ws.Range("I7").Value = 0.001538691255645832
$ws.Range("J7").Value = 0.001538691255645832
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 117.0512696666667
$ws.Range("N7").Value = 351.153809
$ws.Range("O7").Value = 0.2000910950200451
$ws.Range("P7").Value = 0.2000910950200451
$ws.Range("Q7").Value = 4.853764999267667
$ws.Range("R7").Value = 43.683884993409
$ws.Range("S7").Value = 0.0003078784182399427
$ws.Range("T7").Value = 0.0003078784182399427

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Has2"
$ws.Range("C8").Value = "Cd44"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.041467
$ws.Range("H8").Value = 0.124401
$ws.Range("I8").Value = 0.001538691255645832
$ws.Range("J8").Value = 0.001538691255645832
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 171.15883
$ws.Range("N8").Value = 513.47649
$ws.Range("O8").Value = 0.2925842480357353
$ws.Range("P8").Value = 0.2925842480357353
$ws.Range("Q8").Value = 7.097443203609999
$ws.Range("R8").Value = 63.87698883249
$ws.Range("S8").Value = 0.0004501968239922971
$ws.Range("T8").Value = 0.0004501968239922971

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Has2"
$ws.Range("C9").Value = "Cd44"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.041467
$ws.Range("H9").Value = 0.124401
$ws.Range("I9").Value = 0.001538691255645832
$ws.Range("J9").Value = 0.001538691255645832
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 57.695868
$ws.Range("N9").Value = 173.087604
$ws.Range("O9").Value = 0.09862711817759588
$ws.Range("P9").Value = 0.09862711817759588
$ws.Range("Q9").Value = 2.392474558356
$ws.Range("R9").Value = 21.532271025204
$ws.Range("S9").Value = 0.0001517566843094149
$ws.Range("T9").Value = 0.0001517566843094149

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Has2"
$ws.Range("C10").Value = "Cd44"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.363417
$ws.Range("H10").Value = 10.090251
$ws.Range("I10").Value = 0.1248043101017806
$ws.Range("J10").Value = 0.1248043101017806
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 239.0839323333333
$ws.Range("N10").Value = 717.251797
$ws.Range("O10").Value = 0.4086975387666237
$ws.Range("P10").Value = 0.4086975387666237
$ws.Range("Q10").Value = 804.1389624367829
$ws.Range("R10").Value = 7237.250661931046
$ws.Range("S10").Value = 0.05100721436606421
$ws.Range("T10").Value = 0.05100721436606421

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Has2"
$ws.Range("C11").Value = "Cd44"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.363417
$ws.Range("H11").Value = 10.090251
$ws.Range("I11").Value = 0.1248043101017806
$ws.Range("J11").Value = 0.1248043101017806
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 117.0512696666667
$ws.Range("N11").Value = 351.153809
$ws.Range("O11").Value = 0.2000910950200451
$ws.Range("P11").Value = 0.2000910950200451
$ws.Range("Q11").Value = 393.692230268451
$ws.Range("R11").Value = 3543.230072416059
$ws.Range("S11").Value = 0.02497223107148656
$ws.Range("T11").Value = 0.02497223107148656

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Has2"
$ws.Range("C12").Value = "Cd44"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.363417
$ws.Range("H12").Value = 10.090251
$ws.Range("I12").Value = 0.1248043101017806
$ws.Range("J12").Value = 0.1248043101017806
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 171.15883
$ws.Range("N12").Value = 513.47649
$ws.Range("O12").Value = 0.2925842480357353
$ws.Range("P12").Value = 0.2925842480357353
$ws.Range("Q12").Value = 575.6785185221099
$ws.Range("R12").Value = 5181.106666698989
$ws.Range("S12").Value = 0.03651577522274821
$ws.Range("T12").Value = 0.0365157752227482

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Has2"
$ws.Range("C13").Value = "Cd44"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.363417
$ws.Range("H13").Value = 10.090251
$ws.Range("I13").Value = 0.1248043101017806
$ws.Range("J13").Value = 0.1248043101017806
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 57.695868
$ws.Range("N13").Value = 173.087604
$ws.Range("O13").Value = 0.09862711817759588
$ws.Range("P13").Value = 0.09862711817759588
$ws.Range("Q13").Value = 194.055263260956
$ws.Range("R13").Value = 1746.497369348604
$ws.Range("S13").Value = 0.01230908944148164
$ws.Range("T13").Value = 0.01230908944148164
